$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 284 (pushes existing rows 284:349 down to 286:351),
# inheriting formatting (incl. the date style on column D) from the row above.
$ws.Rows("284:285").Insert()

# New "Primera" record (row 284)
$ws.Range("A284").Value = 1
$ws.Range("B284").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C284").Value = "Arica y Parinacota"
$ws.Range("D284").Value = 44785
$ws.Range("E284").Value = 15
$ws.Range("F284").Value = 100114014
$ws.Range("G284").Value = "Betarraga"
$ws.Range("H284").Value = "Sin especificar"
$ws.Range("I284").Value = "Primera"
$ws.Range("J284").Value = 1000
$ws.Range("K284").Value = 500
$ws.Range("L284").Value = 600
$ws.Range("M284").Value = 550
$ws.Range("N284").Value = "`$/paquete 4 unidades"
$ws.Range("O284").Value = "Región de Arica y Parinacota"
$ws.Range("P284").Value = 138
$ws.Range("Q284").Value = 4
$ws.Range("R284").Value = "Hortaliza"

# New "Segunda" record (row 285)
$ws.Range("A285").Value = 1
$ws.Range("B285").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C285").Value = "Arica y Parinacota"
$ws.Range("D285").Value = 44785
$ws.Range("E285").Value = 15
$ws.Range("F285").Value = 100114014
$ws.Range("G285").Value = "Betarraga"
$ws.Range("H285").Value = "Sin especificar"
$ws.Range("I285").Value = "Segunda"
$ws.Range("J285").Value = 1000
$ws.Range("K285").Value = 500
$ws.Range("L285").Value = 600
$ws.Range("M285").Value = 550
$ws.Range("N285").Value = "`$/paquete 5 unidades"
$ws.Range("O285").Value = "Región de Arica y Parinacota"
$ws.Range("P285").Value = 110
$ws.Range("Q285").Value = 5
$ws.Range("R285").Value = "Hortaliza"
